# PS3.xlsx completion pass:
# - fill in Developer/Publisher for L.A. Noire (row 5)
# - append Borderlands: The Pre-Sequel and Grand Theft Auto V
# - insert a fresh blank row at 7 (pushes the Batman: Arkham City block and
#   everything below it down by one) and backfill Developer/Publisher for
#   those shifted rows
# - even up the Developer/Publisher column widths
# - set the page to portrait

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- L.A. Noire gets its Developer/Publisher ---
$ws.Range("D5").Value = "Team Bondi"
$ws.Range("E5").Value = "Rockstar Games"

# --- insert a new row at 7, shifting Batman: Arkham City .. Infamous 2 down ---
$ws.Rows.Item(7).Insert()

# --- two new games appended at the bottom ---
$ws.Range("A15").Value = "Borderlands: The Pre-Sequel"
$ws.Range("E15").Value = "2k Games"
$ws.Range("A16").Value = "Grand Theft Auto V"
$ws.Range("D15").Value = "2K Australia"
$ws.Range("D16").Value = "Rockstar North"
$ws.Range("E16").Value = "Rockstar Games"
$ws.Range("B15").Value = "FPS"
$ws.Range("C15").Value = 2014
$ws.Range("F15").Value = "N"
$ws.Range("B16").Value = "Action"
$ws.Range("C16").Value = 2013

# --- backfill Developer/Publisher on the rows that just shifted down ---
$ws.Range("D8").Value = "Rocksteady Studios"
$ws.Range("E8").Value = "Warner Bros."

$ws.Range("D9").Value = "Monolith"
$ws.Range("E9").Value = "Warner Bros."

$ws.Range("D10").Value = "Visceral Games"
$ws.Range("E10").Value = "EA"

$ws.Range("D11").Value = "Insomniac"
$ws.Range("E11").Value = "Sony Computer"

$ws.Range("D12").Value = "Sucker Punch"
$ws.Range("E12").Value = "Sony Computer"

$ws.Range("D13").Value = "Rocksteady Studios"
$ws.Range("E13").Value = "Warner Bros."

$ws.Range("D14").Value = "Sucker Punch"
$ws.Range("E14").Value = "Sony Computer"

# --- even up the Developer / Publisher column widths ---
$ws.Range("D1:E1").EntireColumn.ColumnWidth = 19.25

# --- leave the cursor where the author left it, and set the sheet to print portrait ---
$ws.Range("F14").Select()
$ws.PageSetup.Orientation = 1
